$d = $word.ActiveDocument

$replacements = @(
    @("76÷7=", "91÷8="),
    @("22÷7=", "47÷3="),
    @("39÷7=", "66÷3="),
    @("12÷4=", "56÷3="),
    @("97÷8=", "68÷8="),
    @("42÷2=", "78÷6="),
    @("83÷8=", "70÷5="),
    @("67÷3=", "93÷9="),
    @("63÷3=", "34÷6="),
    @("25÷9=", "88÷4="),
    @("84÷5=", "70÷9="),
    @("63÷6=", "63÷2="),
    @("52÷7=", "15÷4="),
    @("44÷2=", "70÷4="),
    @("35÷2=", "78÷3="),
    @("49÷6=", "92÷7="),
    @("55÷8=", "65÷4="),
    @("23÷4=", "10÷3="),
    @("26÷5=", "95÷5="),
    @("84÷3=", "51÷7="),
    @("99÷5=", "85÷2="),
    @("19÷8=", "38÷5="),
    @("60÷7=", "61÷5="),
    @("80÷7=", "38÷2="),
    @("31÷3=", "82÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
